$wb = $excel.ActiveWorkbook

# --- Sheet1: "Trends Status" ---
$ws1 = $wb.Worksheets.Item("Trends Status")
$ws1.Range("C7").Value = 10
$ws1.Range("B8").Value = 371
$ws1.Range("C8").Value = 361

# --- Sheet3: "Priority Status" ---
$ws3 = $wb.Worksheets.Item("Priority Status")
$ws3.Range("B2").Value = 103
$ws3.Range("B3").Value = 286
$ws3.Range("B4").Value = 554

# --- Sheet4: "Species qualification" ---
$ws4 = $wb.Worksheets.Item("Species qualification")
$ws4.Range("A2").Value = "SoIB Assessment"
$ws4.Range("B2").Value = 371
$ws4.Range("B4").Value = 10

# --- Sheet5: "High Priority break-up" -> rename + replace data ---
$ws5 = $wb.Worksheets.Item("High Priority break-up")

# Capture the original data (for the new "Major update" sheet) before overwriting.
$origA2 = $ws5.Range("A2").Value2
$origB2 = $ws5.Range("B2").Value2
$origC2 = $ws5.Range("C2").Value2
$origD2 = $ws5.Range("D2").Value2
$origE2 = $ws5.Range("E2").Value2

# Rename sheet5 to "Interannual update - High Pri"
$ws5.Name = "Interannual update - High Pri"

# Clear old data row and write the new "Interannual update" content.
$ws5.Range("A2:E2").ClearContents()

$ws5.Range("A2").Value = "Trend New"
$ws5.Range("B2").Value = 92
$ws5.Range("C2").Value = 89.3
$ws5.Range("D2").Value = 92
$ws5.Range("E2").Value = 93.90000000000001

$ws5.Range("A3").Value = "IUCN"
$ws5.Range("B3").Value = 11
$ws5.Range("C3").Value = 10.7
$ws5.Range("D3").Value = 6
$ws5.Range("E3").Value = 6.1

# --- New Sheet6: "Major update - High Priority " (inserted right after sheet5) ---
$ws6 = $wb.Worksheets.Add($null, $ws5)
$ws6.Name = "Major update - High Priority "

$ws6.Range("A1").Value = "Break-up"
$ws6.Range("B1").Value = "High Species (no.)"
$ws6.Range("C1").Value = "High Species (perc.)"
$ws6.Range("D1").Value = "New High Species (no.)"
$ws6.Range("E1").Value = "New High Species (perc.)"

# Match the bold, centered header style used by every other sheet's row 1.
$ws6.Range("A1:E1").Font.Bold = $true
$ws6.Range("A1:E1").HorizontalAlignment = -4108

$ws6.Range("A2").Value = $origA2
$ws6.Range("B2").Value = $origB2
$ws6.Range("C2").Value = $origC2
$ws6.Range("D2").Value = $origD2
$ws6.Range("E2").Value = $origE2

$wb.Save()
